$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.62"
$ws.Range("D3").Value = "'22.78"
$ws.Range("D4").Value = "'6.165"
$ws.Range("D5").Value = "'0.06085"
$ws.Range("D6").Value = "'6.720"
$ws.Range("D7").Value = "'3.453"
$ws.Range("D8").Value = "'1.354"
$ws.Range("D9").Value = "'0.7970"
$ws.Range("D10").Value = "'0.1580"
$ws.Range("D11").Value = "'0.08035"
$ws.Range("D12").Value = "'0.03344"
$ws.Range("D13").Value = "'0.03088"
$ws.Range("D14").Value = "'0.09299"
$ws.Range("D15").Value = "'3.905"
$ws.Range("D16").Value = "'0.001696"
$ws.Range("D17").Value = "'0.04841"
$ws.Range("D18").Value = "'0.0006150"
$ws.Range("D19").Value = "'0.006221"
$ws.Range("D20").Value = "'0.001101"
$ws.Range("D21").Value = "'0.003382"
$ws.Range("D23").Value = "'3.685"
$ws.Range("D24").Value = "'2.262"
$ws.Range("D26").Value = "'0.1227"
$ws.Range("D27").Value = "'0.0003021"
$ws.Range("D41").Value = "'0.007102"
$ws.Range("D45").Value = "'0.002975"
$ws.Range("D46").Value = "'0.00005987"
$ws.Range("D48").Value = "'0.7513"
$ws.Range("D49").Value = "'0.06599"
